$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = '''15'
$ws.Cells.Item(3, 3).Value = '''2'
$ws.Cells.Item(3, 4).Value = 'https://cdn.sofifa.net/players/020/801/15_120.png'
$ws.Cells.Item(3, 5).Value = 'Cristiano Ronaldo dos Santos Aveiro'
$ws.Cells.Item(3, 6).Value = 'LW, LM'
$ws.Cells.Item(3, 7).Value = 'Real Madrid CF'
$ws.Cells.Item(3, 8).Value = 'Portugal'
$ws.Cells.Item(3, 9).Value = 92
$ws.Cells.Item(3, 10).Value = 92
$ws.Cells.Item(3, 11).Value = 79000000
$ws.Cells.Item(3, 12).Value = 375000
$ws.Cells.Item(3, 13).Value = 29
$ws.Cells.Item(3, 14).Value = 185
$ws.Cells.Item(3, 15).Value = 80
$ws.Cells.Item(3, 16).Value = 'Right'
$ws.Cells.Item(3, 17).Value = 4
$ws.Cells.Item(3, 18).Value = 5
$ws.Cells.Item(3, 19).Value = 5
$ws.Cells.Item(3, 20).Value = 'High/Low'
$ws.Cells.Item(3, 21).Value = 'Normal (185+)'
$ws.Cells.Item(3, 22).Value = 93
$ws.Cells.Item(3, 23).Value = 93
$ws.Cells.Item(3, 24).Value = 81
$ws.Cells.Item(3, 25).Value = 91
$ws.Cells.Item(3, 26).Value = 88
$ws.Cells.Item(3, 27).Value = 79
$ws.Cells.Item(3, 28).Value = 72
$ws.Cells.Item(3, 29).Value = 92
$ws.Cells.Item(3, 30).Value = 91
$ws.Cells.Item(3, 31).Value = 94
$ws.Cells.Item(3, 32).Value = 93
$ws.Cells.Item(3, 33).Value = 90
$ws.Cells.Item(3, 34).Value = 63
$ws.Cells.Item(3, 35).Value = 94
$ws.Cells.Item(3, 36).Value = 94
$ws.Cells.Item(3, 37).Value = 89
$ws.Cells.Item(3, 38).Value = 79
$ws.Cells.Item(3, 39).Value = 93
$ws.Cells.Item(3, 40).Value = 63
$ws.Cells.Item(3, 41).Value = 24
$ws.Cells.Item(3, 42).Value = 91
$ws.Cells.Item(3, 43).Value = 81
$ws.Cells.Item(3, 44).Value = 85
$ws.Cells.Item(3, 45).Value = 0
$ws.Cells.Item(3, 47).Value = 31
$ws.Cells.Item(3, 48).Value = 23
$ws.Cells.Item(3, 49).Value = 7
$ws.Cells.Item(3, 50).Value = 11
$ws.Cells.Item(3, 51).Value = 15
$ws.Cells.Item(3, 52).Value = 14
$ws.Cells.Item(3, 53).Value = 11
$ws.Cells.Item(3, 54).Value = 0
$ws.Cells.Item(3, 55).Value = 'Power Free-Kick, Flair, Long Shot Taker (AI), Speed Dribbler (AI)'

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '''15'
$ws.Cells.Item(4, 3).Value = '''2'
$ws.Cells.Item(4, 4).Value = 'https://cdn.sofifa.net/players/009/014/15_120.png'
$ws.Cells.Item(4, 5).Value = 'Arjen Robben'
$ws.Cells.Item(4, 6).Value = 'RM, LM, RW'
$ws.Cells.Item(4, 7).Value = 'FC Bayern München'
$ws.Cells.Item(4, 8).Value = 'Netherlands'
$ws.Cells.Item(4, 9).Value = 90
$ws.Cells.Item(4, 10).Value = 90
$ws.Cells.Item(4, 11).Value = 54500000
$ws.Cells.Item(4, 12).Value = 275000
$ws.Cells.Item(4, 13).Value = 30
$ws.Cells.Item(4, 14).Value = 180
$ws.Cells.Item(4, 15).Value = 80
$ws.Cells.Item(4, 16).Value = 'Left'
$ws.Cells.Item(4, 17).Value = 2
$ws.Cells.Item(4, 18).Value = 4
$ws.Cells.Item(4, 19).Value = 5
$ws.Cells.Item(4, 20).Value = 'High/Low'
$ws.Cells.Item(4, 21).Value = 'Normal (170-185)'
$ws.Cells.Item(4, 22).Value = 93
$ws.Cells.Item(4, 23).Value = 86
$ws.Cells.Item(4, 24).Value = 83
$ws.Cells.Item(4, 25).Value = 92
$ws.Cells.Item(4, 26).Value = 85
$ws.Cells.Item(4, 27).Value = 83
$ws.Cells.Item(4, 28).Value = 76
$ws.Cells.Item(4, 29).Value = 90
$ws.Cells.Item(4, 30).Value = 93
$ws.Cells.Item(4, 31).Value = 93
$ws.Cells.Item(4, 32).Value = 93
$ws.Cells.Item(4, 33).Value = 89
$ws.Cells.Item(4, 34).Value = 91
$ws.Cells.Item(4, 35).Value = 86
$ws.Cells.Item(4, 36).Value = 61
$ws.Cells.Item(4, 37).Value = 78
$ws.Cells.Item(4, 38).Value = 65
$ws.Cells.Item(4, 39).Value = 90
$ws.Cells.Item(4, 40).Value = 47
$ws.Cells.Item(4, 41).Value = 39
$ws.Cells.Item(4, 42).Value = 89
$ws.Cells.Item(4, 43).Value = 84
$ws.Cells.Item(4, 44).Value = 80
$ws.Cells.Item(4, 45).Value = 0
$ws.Cells.Item(4, 47).Value = 26
$ws.Cells.Item(4, 48).Value = 26
$ws.Cells.Item(4, 49).Value = 10
$ws.Cells.Item(4, 50).Value = 8
$ws.Cells.Item(4, 51).Value = 11
$ws.Cells.Item(4, 52).Value = 5
$ws.Cells.Item(4, 53).Value = 15
$ws.Cells.Item(4, 54).Value = 0
$ws.Cells.Item(4, 55).Value = 'Diver, Injury Prone, Avoids Using Weaker Foot, Selfish, Long Shot Taker (AI), Speed Dribbler (AI), Chip Shot (AI)'

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '''15'
$ws.Cells.Item(5, 3).Value = '''2'
$ws.Cells.Item(5, 4).Value = 'https://cdn.sofifa.net/players/041/236/15_120.png'
$ws.Cells.Item(5, 5).Value = 'Zlatan Ibrahimović'
$ws.Cells.Item(5, 6).Value = 'ST'
$ws.Cells.Item(5, 7).Value = 'Paris Saint-Germain'
$ws.Cells.Item(5, 8).Value = 'Sweden'
$ws.Cells.Item(5, 9).Value = 90
$ws.Cells.Item(5, 10).Value = 90
$ws.Cells.Item(5, 11).Value = 52500000
$ws.Cells.Item(5, 12).Value = 275000
$ws.Cells.Item(5, 13).Value = 32
$ws.Cells.Item(5, 14).Value = 195
$ws.Cells.Item(5, 15).Value = 95
$ws.Cells.Item(5, 16).Value = 'Right'
$ws.Cells.Item(5, 17).Value = 4
$ws.Cells.Item(5, 18).Value = 4
$ws.Cells.Item(5, 19).Value = 5
$ws.Cells.Item(5, 20).Value = 'Medium/Low'
$ws.Cells.Item(5, 21).Value = 'Normal (185+)'
$ws.Cells.Item(5, 22).Value = 76
$ws.Cells.Item(5, 23).Value = 91
$ws.Cells.Item(5, 24).Value = 81
$ws.Cells.Item(5, 25).Value = 86
$ws.Cells.Item(5, 26).Value = 80
$ws.Cells.Item(5, 27).Value = 80
$ws.Cells.Item(5, 28).Value = 76
$ws.Cells.Item(5, 29).Value = 90
$ws.Cells.Item(5, 30).Value = 74
$ws.Cells.Item(5, 31).Value = 77
$ws.Cells.Item(5, 32).Value = 86
$ws.Cells.Item(5, 33).Value = 85
$ws.Cells.Item(5, 34).Value = 41
$ws.Cells.Item(5, 35).Value = 93
$ws.Cells.Item(5, 36).Value = 72
$ws.Cells.Item(5, 37).Value = 78
$ws.Cells.Item(5, 38).Value = 93
$ws.Cells.Item(5, 39).Value = 88
$ws.Cells.Item(5, 40).Value = 84
$ws.Cells.Item(5, 41).Value = 20
$ws.Cells.Item(5, 42).Value = 86
$ws.Cells.Item(5, 43).Value = 83
$ws.Cells.Item(5, 44).Value = 91
$ws.Cells.Item(5, 45).Value = 0
$ws.Cells.Item(5, 47).Value = 41
$ws.Cells.Item(5, 48).Value = 27
$ws.Cells.Item(5, 49).Value = 13
$ws.Cells.Item(5, 50).Value = 15
$ws.Cells.Item(5, 51).Value = 10
$ws.Cells.Item(5, 52).Value = 9
$ws.Cells.Item(5, 53).Value = 12
$ws.Cells.Item(5, 54).Value = 0
$ws.Cells.Item(5, 55).Value = 'Power Free-Kick, Leadership, Flair, Long Shot Taker (AI), Technical Dribbler (AI)'

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = '''15'
$ws.Cells.Item(6, 3).Value = '''2'
$ws.Cells.Item(6, 4).Value = 'https://cdn.sofifa.net/players/167/495/15_120.png'
$ws.Cells.Item(6, 5).Value = 'Manuel Peter Neuer'
$ws.Cells.Item(6, 6).Value = 'GK'
$ws.Cells.Item(6, 7).Value = 'FC Bayern München'
$ws.Cells.Item(6, 8).Value = 'Germany'
$ws.Cells.Item(6, 9).Value = 90
$ws.Cells.Item(6, 10).Value = 90
$ws.Cells.Item(6, 11).Value = 63500000
$ws.Cells.Item(6, 12).Value = 300000
$ws.Cells.Item(6, 13).Value = 28
$ws.Cells.Item(6, 14).Value = 193
$ws.Cells.Item(6, 15).Value = 92
$ws.Cells.Item(6, 16).Value = 'Right'
$ws.Cells.Item(6, 17).Value = 4
$ws.Cells.Item(6, 18).Value = 1
$ws.Cells.Item(6, 19).Value = 5
$ws.Cells.Item(6, 20).Value = 'Medium/Medium'
$ws.Cells.Item(6, 21).Value = 'Normal (185+)'
$ws.Cells.Item(6, 22).Value = 0
$ws.Cells.Item(6, 23).Value = 0
$ws.Cells.Item(6, 24).Value = 0
$ws.Cells.Item(6, 25).Value = 0
$ws.Cells.Item(6, 26).Value = 25
$ws.Cells.Item(6, 27).Value = 25
$ws.Cells.Item(6, 28).Value = 41
$ws.Cells.Item(6, 29).Value = 31
$ws.Cells.Item(6, 30).Value = 58
$ws.Cells.Item(6, 31).Value = 61
$ws.Cells.Item(6, 32).Value = 43
$ws.Cells.Item(6, 33).Value = 89
$ws.Cells.Item(6, 34).Value = 35
$ws.Cells.Item(6, 35).Value = 42
$ws.Cells.Item(6, 36).Value = 78
$ws.Cells.Item(6, 37).Value = 44
$ws.Cells.Item(6, 38).Value = 83
$ws.Cells.Item(6, 39).Value = 25
$ws.Cells.Item(6, 40).Value = 29
$ws.Cells.Item(6, 41).Value = 30
$ws.Cells.Item(6, 42).Value = 25
$ws.Cells.Item(6, 43).Value = 20
$ws.Cells.Item(6, 44).Value = 37
$ws.Cells.Item(6, 45).Value = 0
$ws.Cells.Item(6, 47).Value = 25
$ws.Cells.Item(6, 48).Value = 25
$ws.Cells.Item(6, 49).Value = 87
$ws.Cells.Item(6, 50).Value = 85
$ws.Cells.Item(6, 51).Value = 92
$ws.Cells.Item(6, 52).Value = 90
$ws.Cells.Item(6, 53).Value = 86
$ws.Cells.Item(6, 54).Value = 60
$ws.Cells.Item(6, 55).Value = 'GK Up for Corners, GK Long Throw, 1-on-1 Rush'

Write-Host "rows 3-6 populated"